$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BNoGP")
$ws.Range("A2").Value = "BAU Gas Pumps"
